# Reorder the city rows (rows 4-10, the block between the fixed top two
# rows and the fixed bottom three rows) by Population, descending.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A4:D10")
$sortKey   = $ws.Range("C4:C10")

# xlDescending = 2, xlSortColumns/top-to-bottom = 1, xlNo header = 2
$dataRange.Sort($sortKey, 2, $null, $null, 1, $null, 1, 2)
